# Fill in carrier (D) values for the practice & generic rows, and add the
# new pair_kind (J) / carrier (D) values for the "unique_video" /
# "unique_audio" rows, per the commit "more work towards final product".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows: carrier (D) mirrors the pair's carrier word (K).
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic stimulus rows 6-9 get a pair_kind (J) of "unique_video".
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
# ... and rows 8-9 get "unique_audio".
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New rows 14-17: kind (C) "unique_video" with carrier (D) "can"/"can"/"do"/"do".
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"
$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"
$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"
$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

# New rows 18-21: kind (C) "unique_audio" with carrier (D) "look"/"look"/"where"/"where".
$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"
$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"
$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"
$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
